# Auto-generated Excel COM-interop script
# Applies updated profit/price figures to the Exodus_Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H110").Value = 62827.43
$ws.Range("J110").Value = 62827.43
$ws.Range("L110").Value = 62827.43
$ws.Range("N110").Value = -71007.42999999999
$ws.Range("H123").Value = 66674.55
$ws.Range("J123").Value = 66674.55
$ws.Range("L123").Value = 66674.55
$ws.Range("N123").Value = -76474.55
$ws.Range("H133").Value = 82423.92999999999
$ws.Range("J133").Value = 82423.92999999999
$ws.Range("L133").Value = 82423.92999999999
$ws.Range("N133").Value = -92543.92999999999
$ws.Range("H134").Value = 70068
$ws.Range("J134").Value = 70068
$ws.Range("L134").Value = 70068
$ws.Range("N134").Value = -80208
$ws.Range("H135").Value = 2601
$ws.Range("I135").Value = 2321.4443
$ws.Range("J135").Value = 3104.2
$ws.Range("K135").Value = 20892.9987
$ws.Range("L135").Value = 27937.8
$ws.Range("M135").Value = -18357.9987
$ws.Range("N135").Value = -33007.8
$ws.Range("H136").Value = 63247.75
$ws.Range("J136").Value = 63247.75
$ws.Range("L136").Value = 63247.75
$ws.Range("N136").Value = -73447.75
$ws.Range("H137").Value = 247737.44
$ws.Range("I137").Value = 1901.5834
$ws.Range("K137").Value = 5704.7502
$ws.Range("M137").Value = -3154.7502
$ws.Range("H138").Value = 1662.434
$ws.Range("I138").Value = 1456.186
$ws.Range("J138").Value = 2549.3
$ws.Range("K138").Value = 4368.558
$ws.Range("L138").Value = 7647.900000000001
$ws.Range("M138").Value = 771.442
$ws.Range("N138").Value = -17927.9

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 569.6
$ws.Range("I2").Value = 456.57144
$ws.Range("J2").Value = 833.3333
$ws.Range("K2").Value = 456.57144
$ws.Range("L2").Value = 833.3333
$ws.Range("M2").Value = -343.57144
$ws.Range("N2").Value = -1059.3333
$ws.Range("H44").Value = 90437.39999999999
$ws.Range("I44").Value = 70000
$ws.Range("K44").Value = 70000
$ws.Range("M44").Value = -69512
$ws.Range("H49").Value = 9999
$ws.Range("J49").Value = 9999
$ws.Range("L49").Value = 9999
$ws.Range("N49").Value = -10519
$ws.Range("H52").Value = 62696
$ws.Range("J52").Value = 62696
$ws.Range("L52").Value = 62696
$ws.Range("N52").Value = -63332
$ws.Range("H74").Value = 39608.926
$ws.Range("I74").Value = 68143.8
$ws.Range("K74").Value = 68143.8
$ws.Range("M74").Value = -67269.8
$ws.Range("H77").Value = 39608.926
$ws.Range("I77").Value = 68143.8
$ws.Range("K77").Value = 340719
$ws.Range("M77").Value = -336351
$ws.Range("H104").Value = 35443.4
$ws.Range("J104").Value = 35443.4
$ws.Range("L104").Value = 35443.4
$ws.Range("N104").Value = -42431.4
$ws.Range("H110").Value = 613.9048
$ws.Range("I110").Value = 572.2778
$ws.Range("K110").Value = 572.2778
$ws.Range("M110").Value = 1472.7222
$ws.Range("H116").Value = 569.6
$ws.Range("I116").Value = 456.57144
$ws.Range("J116").Value = 833.3333
$ws.Range("K116").Value = 456.57144
$ws.Range("L116").Value = 833.3333
$ws.Range("M116").Value = 1837.42856
$ws.Range("N116").Value = -5421.3333
$ws.Range("H121").Value = 53797.5
$ws.Range("J121").Value = 53797.5
$ws.Range("L121").Value = 53797.5
$ws.Range("N121").Value = -57291.5
$ws.Range("H122").Value = 4374
$ws.Range("I122").Value = 4513.143
$ws.Range("K122").Value = 13539.429
$ws.Range("M122").Value = -11089.429
$ws.Range("H127").Value = 99707.14
$ws.Range("J127").Value = 99707.14
$ws.Range("L127").Value = 99707.14
$ws.Range("N127").Value = -109627.14

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 30499.334
$ws.Range("J2").Value = 30499.334
$ws.Range("L2").Value = 30499.334
$ws.Range("N2").Value = -30725.334
$ws.Range("H3").Value = 569.6
$ws.Range("I3").Value = 456.57144
$ws.Range("J3").Value = 833.3333
$ws.Range("K3").Value = 456.57144
$ws.Range("L3").Value = 833.3333
$ws.Range("M3").Value = -342.57144
$ws.Range("N3").Value = -1061.3333
$ws.Range("H6").Value = 8250
$ws.Range("I6").Value = 10000
$ws.Range("J6").Value = 6500
$ws.Range("K6").Value = 10000
$ws.Range("L6").Value = 6500
$ws.Range("M6").Value = -9887
$ws.Range("N6").Value = -6726
$ws.Range("H50").Value = 62675.25
$ws.Range("J50").Value = 66664
$ws.Range("L50").Value = 66664
$ws.Range("N50").Value = -67812
$ws.Range("H51").Value = 34997.332
$ws.Range("J51").Value = 34997.332
$ws.Range("L51").Value = 34997.332
$ws.Range("N51").Value = -35979.332
$ws.Range("H52").Value = 99988
$ws.Range("J52").Value = 99988
$ws.Range("L52").Value = 99988
$ws.Range("N52").Value = -100514
$ws.Range("H109").Value = 88327.5
$ws.Range("J109").Value = 88327.5
$ws.Range("L109").Value = 88327.5
$ws.Range("N109").Value = -91101.5
$ws.Range("H119").Value = 99986
$ws.Range("J119").Value = 99986
$ws.Range("L119").Value = 99986
$ws.Range("N119").Value = -109662
$ws.Range("H121").Value = 99988
$ws.Range("J121").Value = 99988
$ws.Range("L121").Value = 99988
$ws.Range("N121").Value = -103482
$ws.Range("H132").Value = 52059.695
$ws.Range("J132").Value = 52059.695
$ws.Range("L132").Value = 52059.695
$ws.Range("N132").Value = -62179.695
$ws.Range("H135").Value = 74710.8
$ws.Range("J135").Value = 74710.8
$ws.Range("L135").Value = 74710.8
$ws.Range("N135").Value = -84850.8
$ws.Range("H138").Value = 89991.336
$ws.Range("J138").Value = 89991.336
$ws.Range("L138").Value = 89991.336
$ws.Range("N138").Value = -100271.336

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 31813
$ws.Range("J9").Value = 31813
$ws.Range("L9").Value = 31813
$ws.Range("N9").Value = -32149
$ws.Range("H22").Value = 987
$ws.Range("I22").Value = 1085.2858
$ws.Range("K22").Value = 1085.2858
$ws.Range("M22").Value = -735.2858000000001
$ws.Range("H58").Value = 1859.7
$ws.Range("I58").Value = 1730.9166
$ws.Range("K58").Value = 1730.9166
$ws.Range("M58").Value = -1527.9166
$ws.Range("H108").Value = 56014.91
$ws.Range("J108").Value = 56014.91
$ws.Range("L108").Value = 56014.91
$ws.Range("N108").Value = -63694.91
$ws.Range("H116").Value = 97735.60000000001
$ws.Range("J116").Value = 97735.60000000001
$ws.Range("L116").Value = 97735.60000000001
$ws.Range("N116").Value = -106913.6
$ws.Range("H136").Value = 1859.7
$ws.Range("I136").Value = 1730.9166
$ws.Range("K136").Value = 5192.7498
$ws.Range("M136").Value = -2642.7498

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 466.66666
$ws.Range("J107").Value = 466.66666
$ws.Range("L107").Value = 1399.99998
$ws.Range("N107").Value = -5239.999980000001
$ws.Range("H113").Value = 67610.87
$ws.Range("I113").Value = 831.4
$ws.Range("J113").Value = 201169.8
$ws.Range("K113").Value = 2494.2
$ws.Range("L113").Value = 603509.3999999999
$ws.Range("M113").Value = -324.1999999999998
$ws.Range("N113").Value = -607849.3999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H56").Value = 49995
$ws.Range("J56").Value = 49995
$ws.Range("L56").Value = 49995
$ws.Range("N56").Value = -51499
$ws.Range("H107").Value = 1022.2857
$ws.Range("I107").Value = 970.2222
$ws.Range("J107").Value = 1061.3334
$ws.Range("K107").Value = 970.2222
$ws.Range("L107").Value = 1061.3334
$ws.Range("M107").Value = 949.7778
$ws.Range("N107").Value = -4901.3334
$ws.Range("H135").Value = 67979.27
$ws.Range("J135").Value = 67979.27
$ws.Range("L135").Value = 67979.27
$ws.Range("N135").Value = -78119.27
$ws.Range("H140").Value = 98496
$ws.Range("J140").Value = 98496
$ws.Range("L140").Value = 98496
$ws.Range("N140").Value = -108856

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H118").Value = 53960.4
$ws.Range("J118").Value = 53960.4
$ws.Range("L118").Value = 53960.4
$ws.Range("N118").Value = -57274.4
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 10000
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 10000
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 10000
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -11214
$ws.Range("H58").Value = 23729.2
$ws.Range("J58").Value = 28346.5
$ws.Range("L58").Value = 28346.5
$ws.Range("N58").Value = -28962.5
